# directed_trips_regions_bimonthly_HCR_plus1.xlsx
# "adjusted projected harvest per trip, tested model"
#
# The projected-harvest-per-trip column (F) used to be a plain R+1
# calculation. It is now rounded to the nearest whole trip:
#   =ROUND(R<row>+1, 0)
# instead of
#   =R<row>+1
#
# Apply as three fill-down blocks (matching the shared-formula runs the
# workbook already used) so the formula text + calculated values line up
# row-for-row with column R.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F13").Formula = "=ROUND(R2+1, 0)"
$ws.Range("F14").Formula = "=ROUND(R14+1, 0)"
$ws.Range("F15:F78").Formula = "=ROUND(R15+1, 0)"
$ws.Range("F79:F124").Formula = "=ROUND(R79+1, 0)"

# Restore the selection used while reviewing the updated formulas
# (previously scrolled down to row 94 selecting the whole F column).
$ws.Range("F2:F14").Select() | Out-Null
